# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Each row is a crypto-currency entry; D = Price (text), E = Volume(1h) % (text).
# Numeric-looking text values are prefixed with a leading apostrophe so Excel
# keeps them as text (matching the original inlineStr/text cell type) instead of
# coercing them to a Number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "45.059.62"
$ws.Range("E2").Value = "  -3.63%  "

# Row 3
$ws.Range("D3").Value = "2.430.09"
$ws.Range("E3").Value = "  +7.43%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'293.86"
$ws.Range("E5").Value = "  -2.27%  "

# Row 6
$ws.Range("D6").Value = "'92.96"
$ws.Range("E6").Value = "  -7.28%  "

# Row 7
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("E9").Value = "  -2.04%  "

# Row 10
$ws.Range("E10").Value = "  -3.85%  "

# Row 11
$ws.Range("E11").Value = "  -0.35%  "

# Row 12
$ws.Range("D12").Value = "'7.01"
$ws.Range("E12").Value = "  -2.00%  "

# Row 13
$ws.Range("E13").Value = "  +1.84%  "

# Row 14
$ws.Range("D14").Value = "2.801.59"
$ws.Range("E14").Value = "  +7.46%  "

# Row 15
$ws.Range("D15").Value = "2.417.91"
$ws.Range("E15").Value = "  +6.71%  "

# Row 16
$ws.Range("D16").Value = "'14.28"
$ws.Range("E16").Value = "  +5.21%  "

# Row 17
$ws.Range("E17").Value = "  +5.51%  "

# Row 18
$ws.Range("D18").Value = "45.060.84"
$ws.Range("E18").Value = "  -3.64%  "

# Row 19
$ws.Range("D19").Value = "'12.37"
$ws.Range("E19").Value = "  -2.52%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0937"

# Row 21
$ws.Range("D21").Value = "'6.21"
$ws.Range("E21").Value = "  +6.08%  "

# Row 22
$ws.Range("D22").Value = "'67.10"
$ws.Range("E22").Value = "  +3.08%  "

# Row 23
$ws.Range("D23").Value = "'239.32"
$ws.Range("E23").Value = "  -3.85%  "

# Row 24
$ws.Range("E24").Value = "  -1.61%  "

# Row 25
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.09%  "

# Row 26
$ws.Range("D26").Value = "'1.92"
$ws.Range("E26").Value = "  +2.62%  "

# Row 27
$ws.Range("E27").Value = "  -0.64%  "

# Row 28
$ws.Range("D28").Value = "'37.16"
$ws.Range("E28").Value = "  -12.62%  "

# Row 29
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("D30").Value = "'3.83"
$ws.Range("E30").Value = "  +20.33%  "

# Row 31
$ws.Range("D31").Value = "'21.33"
$ws.Range("E31").Value = "  +7.61%  "

# Row 32
$ws.Range("D32").Value = "'148.80"
$ws.Range("E32").Value = "  +2.56%  "

# Row 33
$ws.Range("E33").Value = "  -2.42%  "

# Row 34
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0761"
$ws.Range("E35").Value = "  -1.31%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.99"
$ws.Range("E36").Value = "  +16.97%  "

# Row 37
$ws.Range("E37").Value = "  -1.86%  "

# Row 38
$ws.Range("E38").Value = "  -0.42%  "

# Row 39
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'14.18"
$ws.Range("E39").Value = "  -11.83%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'3.73"
$ws.Range("E40").Value = "  -2.49%  "

# Row 41
$ws.Range("E41").Value = "  -1.24%  "

# Row 42
$ws.Range("D42").Value = "1.991.77"
$ws.Range("E42").Value = "  +11.21%  "

# Row 43
$ws.Range("D43").Value = "'3.17"
$ws.Range("E43").Value = "  -1.11%  "

# Row 44
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("D45").Value = "'88.48"
$ws.Range("E45").Value = "  -2.67%  "

# Row 46
$ws.Range("D46").Value = "'16.07"
$ws.Range("E46").Value = "  +26.30%  "

# Row 47
$ws.Range("E47").Value = "  -13.67%  "

# Row 48
$ws.Range("D48").Value = "'8.57"
$ws.Range("E48").Value = "  +10.03%  "

# Row 49
$ws.Range("D49").Value = "'101.35"
$ws.Range("E49").Value = "  +8.37%  "

# Row 50
$ws.Range("D50").Value = "2.671.50"
$ws.Range("E50").Value = "  +7.47%  "

# Row 51
$ws.Range("D51").Value = "'0.181"
$ws.Range("E51").Value = "  -3.57%  "
